# Sprite designer "DIGGING" animation frames update.
# The worksheet encodes 3 small 8x9 bitmaps (rows 11-19, 21-29, 31-39,
# columns A-H) as 1/blank toggle cells; J:M / O columns are formulas that
# derive byte values from those bitmaps and recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$cellsToSet = @(
    "F12","B15","C15","F15","A16","D16","G16","A17","C17","E17","G17",
    "C18","E18","B19","C19",
    "F22","B25","C25","F25","A26","G26","A27","C27","G27","C28","E28","B29","C29",
    "F32","B35","C35","F35","A36","G36","A37","G37","B38","E38","E39","F39"
)

$cellsToClear = @(
    "E12","B16","F17","F18",
    "E22","B26","B27","F27","A28","B28","F29",
    "E32","B39","C39"
)

foreach ($ref in $cellsToSet) {
    $ws.Range($ref).Value = 1
}

foreach ($ref in $cellsToClear) {
    $ws.Range($ref).Value = ""
}

# Restore the selection Excel records after editing the third sprite block.
$ws.Range("O31:O33").Select()
